# Apply edits described by the commit diff: populate the "G" (quantity) column
# for several billing line items on Sheet1. The dependent "I" column formulas
# will recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G9").Value  = 117
$ws.Range("G12").Value = 117
$ws.Range("G14").Value = 119
$ws.Range("G16").Value = 27
$ws.Range("G17").Value = 25
$ws.Range("G18").Value = 118
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1
$ws.Range("G27").Value = 4
$ws.Range("G29").Value = 10

$excel.CalculateFull()
